# Update quarterly report: drop oldest quarter (1399/06) and append the
# newest quarter (1401/12), shifting all quarterly figures one column to
# the left (E:M <- F:N) and filling in the freshly computed figures for
# the new quarter / recalculated figures from the updated read_price
# algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$cols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N")

$quarters = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

$row10 = @(10783849, 8399574, 8033579, 9452568, 4119370, 7244035, 11613330, 11715998, 7839969, 8192960)
$row14 = @(10872, 15478, 9159, 20283, 24775, 36081, 10399, 42413, 5617, 27139)
$row16 = @(4847, 155073, 89032, 126746, 114221, -13328, 76272, 76272, 76271, 929476)
$row17 = @(103911, 150168, 200712, 261075, 77802, 290192, 312800, 977844, -365553, 511652)
$row19 = @(614130, 629700, 519399, 804525, 423321, 984355, 814377, 954188, 813405, 1676192)
$row20 = @(11517609, 9349993, 8851881, 10665197, 4759489, 8541335, 12827178, 13766715, 8369709, 11337419)
$row26 = @(701, 702, 701, 692, 697, 709, 709, 710, 710, 707)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $c = $cols[$i]
    $ws.Range($c + "8").Value = $quarters[$i]
    $ws.Range($c + "24").Value = $quarters[$i]
    $ws.Range($c + "10").Value = $row10[$i]
    $ws.Range($c + "14").Value = $row14[$i]
    $ws.Range($c + "16").Value = $row16[$i]
    $ws.Range($c + "17").Value = $row17[$i]
    $ws.Range($c + "19").Value = $row19[$i]
    $ws.Range($c + "20").Value = $row20[$i]
    $ws.Range($c + "26").Value = $row26[$i]
}
